# Weekly data refresh: insert two new rows of "Apio" data at the top of the
# existing data block (rows 413-414), pushing the previous rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 413; this shifts the existing rows 413-438
# down to 415-440 and carries formatting (incl. the date-format style on
# column D) from the row above into the new blank rows.
$ws.Rows("413:414").Insert()

# Populate new row 413.
$ws.Range("A413").Value = 11
$ws.Range("B413").Value = "Vega Monumental Concepción"
$ws.Range("C413").Value = "Bíobío"
$ws.Range("D413").Value = 45013
$ws.Range("E413").Value = 8
$ws.Range("F413").Value = 100112017
$ws.Range("G413").Value = "Apio"
$ws.Range("H413").Value = "Americana (o)"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 220
$ws.Range("K413").Value = 7000
$ws.Range("L413").Value = 8000
$ws.Range("M413").Value = 7545
$ws.Range("N413").Value = "$/docena de matas"
$ws.Range("O413").Value = "Región de Coquimbo"
$ws.Range("P413").Value = 1258
$ws.Range("Q413").Value = 6
$ws.Range("R413").Value = "Hortaliza"

# Populate new row 414.
$ws.Range("A414").Value = 11
$ws.Range("B414").Value = "Vega Monumental Concepción"
$ws.Range("C414").Value = "Bíobío"
$ws.Range("D414").Value = 45013
$ws.Range("E414").Value = 8
$ws.Range("F414").Value = 100112017
$ws.Range("G414").Value = "Apio"
$ws.Range("H414").Value = "Americana (o)"
$ws.Range("I414").Value = "Segunda"
$ws.Range("J414").Value = 150
$ws.Range("K414").Value = 5000
$ws.Range("L414").Value = 5000
$ws.Range("M414").Value = 5000
$ws.Range("N414").Value = "$/docena de matas"
$ws.Range("O414").Value = "Región de Coquimbo"
$ws.Range("P414").Value = 833
$ws.Range("Q414").Value = 6
$ws.Range("R414").Value = "Hortaliza"
